$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 798.6667
$ws.Range("I92").Value = 886.2857
$ws.Range("J92").Value = 676
$ws.Range("K92").Value = 886.2857
$ws.Range("L92").Value = 676
$ws.Range("M92").Value = 361.7143
$ws.Range("N92").Value = -3172
# Row 98
$ws.Range("H98").Value = 3983
$ws.Range("I98").Value = 2021.091
$ws.Range("J98").Value = 8299.200000000001
$ws.Range("K98").Value = 2021.091
$ws.Range("L98").Value = 8299.200000000001
$ws.Range("M98").Value = -523.0909999999999
$ws.Range("N98").Value = -11295.2
# Row 122
$ws.Range("H122").Value = 3983
$ws.Range("I122").Value = 2021.091
$ws.Range("J122").Value = 8299.200000000001
$ws.Range("K122").Value = 6063.272999999999
$ws.Range("L122").Value = 24897.6
$ws.Range("M122").Value = -3613.272999999999
$ws.Range("N122").Value = -29797.6
# Row 129
$ws.Range("H129").Value = 1043.1951
$ws.Range("I129").Value = 897
$ws.Range("J129").Value = 1046.85
$ws.Range("K129").Value = 2691
$ws.Range("L129").Value = 3140.55
$ws.Range("M129").Value = 2309
$ws.Range("N129").Value = -13140.55
# Row 133
$ws.Range("H133").Value = 47217.777
$ws.Range("J133").Value = 47217.777
$ws.Range("L133").Value = 47217.777
$ws.Range("N133").Value = -57337.777

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 863.73334
$ws.Range("I97").Value = 862.4828
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 862.4828
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -366.4828
$ws.Range("N97").Value = -1892
# Row 139
$ws.Range("H139").Value = 83000
$ws.Range("J139").Value = 83000
$ws.Range("L139").Value = 83000
$ws.Range("N139").Value = -93280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1624.1666
$ws.Range("I86").Value = 1549.1364
$ws.Range("J86").Value = 1830.5
$ws.Range("K86").Value = 1549.1364
$ws.Range("L86").Value = 1830.5
$ws.Range("M86").Value = -426.1364000000001
$ws.Range("N86").Value = -4076.5
# Row 89
$ws.Range("H89").Value = 1624.1666
$ws.Range("I89").Value = 1549.1364
$ws.Range("J89").Value = 1830.5
$ws.Range("K89").Value = 7745.682000000001
$ws.Range("L89").Value = 9152.5
$ws.Range("M89").Value = -2129.682000000001
$ws.Range("N89").Value = -20384.5
# Row 107
$ws.Range("H107").Value = 1568.8334
$ws.Range("I107").Value = 1382.6
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 1382.6
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 537.4000000000001
$ws.Range("N107").Value = -6340

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 292.125
$ws.Range("J22").Value = 962.6
$ws.Range("K22").Value = 292.125
$ws.Range("L22").Value = 962.6
$ws.Range("M22").Value = 57.875
$ws.Range("N22").Value = -1662.6
# Row 70
$ws.Range("H70").Value = 27995
$ws.Range("J70").Value = 27995
$ws.Range("L70").Value = 27995
$ws.Range("N70").Value = -28625
# Row 73
$ws.Range("H73").Value = 27995
$ws.Range("J73").Value = 27995
$ws.Range("L73").Value = 27995
$ws.Range("N73").Value = -30179

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 126.63158
$ws.Range("I23").Value = 92.57143000000001
$ws.Range("J23").Value = 146.5
$ws.Range("K23").Value = 277.71429
$ws.Range("L23").Value = 439.5
$ws.Range("M23").Value = -42.71429000000001
$ws.Range("N23").Value = -909.5
# Row 107
$ws.Range("H107").Value = 844.2857
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
# Row 113
$ws.Range("H113").Value = 1010664.7
$ws.Range("I113").Value = 1894500.4
$ws.Range("J113").Value = 566.7857
$ws.Range("K113").Value = 5683501.199999999
$ws.Range("L113").Value = 1700.3571
$ws.Range("M113").Value = -5681331.199999999
$ws.Range("N113").Value = -6040.3571
# Row 126
$ws.Range("H126").Value = 5699.125
$ws.Range("I126").Value = 2976.6667
$ws.Range("J126").Value = 7332.6
$ws.Range("K126").Value = 8930.000100000001
$ws.Range("L126").Value = 21997.8
$ws.Range("M126").Value = -3990.000100000001
$ws.Range("N126").Value = -31877.8
# Row 137
$ws.Range("H137").Value = 12265062
$ws.Range("I137").Value = 3054.6155
$ws.Range("J137").Value = 22892136
$ws.Range("K137").Value = 9163.8465
$ws.Range("L137").Value = 68676408
$ws.Range("M137").Value = -4063.8465
$ws.Range("N137").Value = -68686608

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 125.454544
$ws.Range("I2").Value = 46.875
$ws.Range("K2").Value = 46.875
$ws.Range("M2").Value = 66.125
# Row 80
$ws.Range("H80").Value = 2204.16
$ws.Range("I80").Value = 2208.7
$ws.Range("J80").Value = 2201.1333
$ws.Range("K80").Value = 2208.7
$ws.Range("L80").Value = 2201.1333
$ws.Range("M80").Value = -1210.7
$ws.Range("N80").Value = -4197.1333
# Row 83
$ws.Range("H83").Value = 2204.16
$ws.Range("I83").Value = 2208.7
$ws.Range("J83").Value = 2201.1333
$ws.Range("K83").Value = 11043.5
$ws.Range("L83").Value = 11005.6665
$ws.Range("M83").Value = -6051.5
$ws.Range("N83").Value = -20989.6665
# Row 107
$ws.Range("H107").Value = 1009.9231
$ws.Range("I107").Value = 605
$ws.Range("J107").Value = 1357
$ws.Range("K107").Value = 605
$ws.Range("L107").Value = 1357
$ws.Range("M107").Value = 1315
$ws.Range("N107").Value = -5197

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 778.3158
$ws.Range("I22").Value = 694.5
$ws.Range("J22").Value = 800.6667
$ws.Range("K22").Value = 694.5
$ws.Range("L22").Value = 800.6667
$ws.Range("M22").Value = -399.5
$ws.Range("N22").Value = -1390.6667
# Row 27
$ws.Range("H27").Value = 778.3158
$ws.Range("I27").Value = 694.5
$ws.Range("J27").Value = 800.6667
$ws.Range("K27").Value = 694.5
$ws.Range("L27").Value = 800.6667
$ws.Range("M27").Value = -587.5
$ws.Range("N27").Value = -1014.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 29766.666
$ws.Range("I75").Value = 29500
$ws.Range("J75").Value = 29900
$ws.Range("K75").Value = 29500
$ws.Range("L75").Value = 29900
$ws.Range("M75").Value = -28564
$ws.Range("N75").Value = -31772
# Row 78
$ws.Range("H78").Value = 29766.666
$ws.Range("I78").Value = 29500
$ws.Range("J78").Value = 29900
$ws.Range("K78").Value = 88500
$ws.Range("L78").Value = 89700
$ws.Range("M78").Value = -83820
$ws.Range("N78").Value = -99060
# Row 132
$ws.Range("H132").Value = 906972.4
$ws.Range("I132").Value = 1280535.6
$ws.Range("J132").Value = 2555.9473
$ws.Range("K132").Value = 3841606.8
$ws.Range("L132").Value = 7667.841899999999
$ws.Range("M132").Value = -3839076.8
$ws.Range("N132").Value = -12727.8419
# Row 136
$ws.Range("H136").Value = 6878.1113
$ws.Range("I136").Value = 3437.0527
$ws.Range("J136").Value = 15050.625
$ws.Range("K136").Value = 10311.1581
$ws.Range("L136").Value = 45151.875
$ws.Range("M136").Value = -7761.158100000001
$ws.Range("N136").Value = -50251.875
